$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 4375
$ws.Range("J29").Value = 5800
$ws.Range("L29").Value = 17400
$ws.Range("N29").Value = -17962
$ws.Range("H55").Value = 307.85715
$ws.Range("J55").Value = 519.6667
$ws.Range("L55").Value = 519.6667
$ws.Range("N55").Value = -947.6667
$ws.Range("H107").Value = 845.2
$ws.Range("I107").Value = 833.7143
$ws.Range("K107").Value = 833.7143
$ws.Range("M107").Value = 1086.2857
$ws.Range("H125").Value = 949.2727
$ws.Range("I125").Value = 930.625
$ws.Range("J125").Value = 999
$ws.Range("K125").Value = 8375.625
$ws.Range("L125").Value = 8991
$ws.Range("M125").Value = -5915.625
$ws.Range("N125").Value = -13911
$ws.Range("H138").Value = 859.25
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1226.52
$ws.Range("I2").Value = 1138.4375
$ws.Range("K2").Value = 1138.4375
$ws.Range("M2").Value = -1025.4375
$ws.Range("H32").Value = 5974
$ws.Range("I32").Value = 6426.1816
$ws.Range("K32").Value = 6426.1816
$ws.Range("M32").Value = -6139.1816
$ws.Range("H68").Value = 44666.668
$ws.Range("J68").Value = 44666.668
$ws.Range("L68").Value = 44666.668
$ws.Range("H71").Value = 44666.668
$ws.Range("J71").Value = 44666.668
$ws.Range("L71").Value = 134000.004
$ws.Range("H116").Value = 1226.52
$ws.Range("I116").Value = 1138.4375
$ws.Range("K116").Value = 1138.4375
$ws.Range("M116").Value = 1155.5625
$ws.Range("H122").Value = 3040.5417
$ws.Range("I122").Value = 3223
$ws.Range("J122").Value = 2128.25
$ws.Range("K122").Value = 9669
$ws.Range("L122").Value = 6384.75
$ws.Range("M122").Value = -7219
$ws.Range("N122").Value = -11284.75
$ws.Range("H132").Value = 3544.5454
$ws.Range("I132").Value = 4073.2307
$ws.Range("K132").Value = 12219.6921
$ws.Range("M132").Value = -9689.6921
$ws.Range("N68").Value = -46288.668
$ws.Range("N71").Value = -142112.004

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1226.52
$ws.Range("I3").Value = 1138.4375
$ws.Range("K3").Value = 1138.4375
$ws.Range("M3").Value = -1024.4375
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1890.1818
$ws.Range("I16").Value = 1782.1666
$ws.Range("K16").Value = 1782.1666
$ws.Range("M16").Value = -1495.1666
$ws.Range("H22").Value = 294.33334
$ws.Range("I22").Value = 294.33334
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 294.33334
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 55.66665999999998
$ws.Range("N22").ClearContents()
$ws.Range("H31").Value = 2663256.8
$ws.Range("I31").Value = 2419.1143
$ws.Range("J31").Value = 10424033
$ws.Range("K31").Value = 2419.1143
$ws.Range("L31").Value = 10424033
$ws.Range("M31").Value = -2124.1143
$ws.Range("N31").Value = -10424623
$ws.Range("H34").Value = 2663256.8
$ws.Range("I34").Value = 2419.1143
$ws.Range("J34").Value = 10424033
$ws.Range("K34").Value = 2419.1143
$ws.Range("L34").Value = 10424033
$ws.Range("M34").Value = -2217.1143
$ws.Range("N34").Value = -10424437
$ws.Range("H107").Value = 3572599.5
$ws.Range("I107").Value = 5556366
$ws.Range("K107").Value = 5556366
$ws.Range("M107").Value = -5554446
$ws.Range("H113").Value = 1890.1818
$ws.Range("I113").Value = 1782.1666
$ws.Range("K113").Value = 1782.1666
$ws.Range("M113").Value = 387.8334
$ws.Range("H122").Value = 388.45
$ws.Range("I122").Value = 343.73334
$ws.Range("K122").Value = 1031.20002
$ws.Range("M122").Value = 1418.79998
$ws.Range("H134").Value = 4495.0835
$ws.Range("I134").Value = 4798.1055
$ws.Range("J134").Value = 3343.6
$ws.Range("K134").Value = 14394.3165
$ws.Range("L134").Value = 10030.8
$ws.Range("M134").Value = -11859.3165
$ws.Range("N134").Value = -15100.8
$ws.Range("H138").Value = 69998.586
$ws.Range("J138").Value = 69998.586
$ws.Range("L138").Value = 69998.586
$ws.Range("N138").Value = -80278.586
$ws.Range("H139").Value = 41999.5
$ws.Range("J139").Value = 41999.5
$ws.Range("L139").Value = 41999.5
$ws.Range("H140").Value = 69908.17999999999
$ws.Range("N139").Value = -52279.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1756.2858
$ws.Range("I75").Value = 1400
$ws.Range("J75").Value = 1898.8
$ws.Range("K75").Value = 4200
$ws.Range("L75").Value = 5696.4
$ws.Range("M75").Value = -3202
$ws.Range("N75").Value = -7692.4
$ws.Range("H78").Value = 1756.2858
$ws.Range("I78").Value = 1400
$ws.Range("J78").Value = 1898.8
$ws.Range("K78").Value = 12600
$ws.Range("L78").Value = 17089.2
$ws.Range("M78").Value = -7608
$ws.Range("N78").Value = -27073.2
$ws.Range("H136").Value = 1200.5
$ws.Range("I136").Value = 1200.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3601.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 1498.5
$ws.Range("N136").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 10269.1
$ws.Range("J107").Value = 14828.5
$ws.Range("L107").Value = 14828.5
$ws.Range("N107").Value = -18668.5
$ws.Range("H113").Value = 2403.625
$ws.Range("I113").Value = 2372
$ws.Range("J113").Value = 2498.5
$ws.Range("K113").Value = 2372
$ws.Range("L113").Value = 2498.5
$ws.Range("M113").Value = -202
$ws.Range("N113").Value = -6838.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2755.0588
$ws.Range("I100").Value = 2250.077
$ws.Range("J100").Value = 4396.25
$ws.Range("K100").Value = 2250.077
$ws.Range("L100").Value = 4396.25
$ws.Range("M100").Value = -1709.077
$ws.Range("N100").Value = -5478.25
$ws.Range("H122").Value = 6297.9287
$ws.Range("J122").Value = 7429.125
$ws.Range("L122").Value = 22287.375
$ws.Range("N122").Value = -27187.375
$ws.Range("H132").Value = 4772.303
$ws.Range("J132").Value = 5375.722
$ws.Range("L132").Value = 16127.166
$ws.Range("N132").Value = -21187.166
$ws.Range("H136").Value = 4599
$ws.Range("I136").Value = 2331.6667
$ws.Range("K136").Value = 6995.000100000001
$ws.Range("M136").Value = -4445.000100000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 37449.5
$ws.Range("J63").Value = 37449.5
$ws.Range("L63").Value = 37449.5
$ws.Range("N63").Value = -38697.5
$ws.Range("H66").Value = 37449.5
$ws.Range("J66").Value = 37449.5
$ws.Range("L66").Value = 112348.5
$ws.Range("N66").Value = -118588.5
$ws.Range("H81").Value = 6775.8887
$ws.Range("I81").Value = 10248.25
$ws.Range("K81").Value = 20496.5
$ws.Range("M81").Value = -19435.5
$ws.Range("H84").Value = 6775.8887
$ws.Range("I84").Value = 10248.25
$ws.Range("K84").Value = 102482.5
$ws.Range("M84").Value = -97178.5
$ws.Range("H100").Value = 71429210
$ws.Range("I100").Value = 738.2222
$ws.Range("K100").Value = 1476.4444
$ws.Range("M100").Value = -935.4444000000001
$ws.Range("H122").Value = 27780620
$ws.Range("I122").Value = 3180.6667
$ws.Range("J122").Value = 83335500
$ws.Range("K122").Value = 9542.000100000001
$ws.Range("L122").Value = 250006500
$ws.Range("M122").Value = -7092.000100000001
$ws.Range("N122").Value = -250011400
$ws.Range("H126").Value = 1820.8
$ws.Range("I126").Value = 2166.3333
$ws.Range("K126").Value = 6498.999899999999
$ws.Range("M126").Value = -4028.999899999999
$ws.Range("H132").Value = 1249.1842
$ws.Range("I132").Value = 1146.9667
$ws.Range("K132").Value = 3440.9001
$ws.Range("M132").Value = -910.9000999999998
$ws.Range("H133").Value = 107999
$ws.Range("J133").Value = 107999
$ws.Range("L133").Value = 107999
$ws.Range("N133").Value = -118119
